$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "Recorded By" contributor lists (G column) ---
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G3").Value = "hend_mahmoud@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G6").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G7").Value = "Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G30").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"

# --- Class Statistics block (K/L) ---
$ws.Range("L6").Value = 27
$ws.Range("L8").Value = 0

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "93.1%"

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "26.9%"

# --- Group Statistics row (row 15) ---
$ws.Range("O15").Value = 27
$ws.Range("Q15").Value = 0

$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "93.1%"

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "26.9%"

# --- Row 19: PARASITOLOGY session 7 is now recorded ---
# Pick up the row-2 (green "Recorded") formatting for the whole row.
$src = $ws.Range("A2:I2")
$dst = $ws.Range("A19:I19")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G19").Value = "Rania.a.youssef@med.asu.edu.eg"
$ws.Range("H19").Value = "18/251"
$ws.Range("I19").Value = "Recorded"
